$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3579947948455811
$ws.Range("E2").Value = 65.65520729320997
$ws.Range("F2").Value = 0.00224943021628463
$ws.Range("G2").Value = 0.001826715171882966
$ws.Range("H2").Value = 0.001826715171882966
$ws.Range("I2").Value = 0.00170740620058687
$ws.Range("J2").Value = 0.001643657613261553
$ws.Range("K2").Value = 0.001421829698083137
$ws.Range("L2").Value = 0.001421829698083137
$ws.Range("M2").Value = 0.001421829698083137
$ws.Range("N2").Value = 0.001421829698083137
$ws.Range("O2").Value = 0.001407940938066215
$ws.Range("P2").Value = 0.001407940938066215
$ws.Range("Q2").Value = 0.001365617171863614
$ws.Range("R2").Value = 0.001348345802293301
$ws.Range("S2").Value = 0.001348345802293301
$ws.Range("T2").Value = 0.001324707258688762
$ws.Range("U2").Value = 0.001324707258688762
$ws.Range("V2").Value = 0.001319506892207172
$ws.Range("W2").Value = 0.001307789499457221
$ws.Range("X2").Value = 0.001304440380025832
$ws.Range("Y2").Value = 0.001279828602206822

$ws.Range("C3").Value = 0.4269788265228271
$ws.Range("E3").Value = 59.17800219631863
$ws.Range("F3").Value = 0.002135335730840912
$ws.Range("G3").Value = 0.001772538797174516
$ws.Range("H3").Value = 0.001729328094364473
$ws.Range("I3").Value = 0.001636737916744652
$ws.Range("J3").Value = 0.001559823603102524
$ws.Range("K3").Value = 0.001559823603102524
$ws.Range("L3").Value = 0.001501951401846489
$ws.Range("M3").Value = 0.00145731077752669
$ws.Range("N3").Value = 0.001301151539918877
$ws.Range("O3").Value = 0.001301151539918877
$ws.Range("P3").Value = 0.001301151539918877
$ws.Range("Q3").Value = 0.001247742550202562
$ws.Range("R3").Value = 0.001247742550202562
$ws.Range("S3").Value = 0.00122737665481305
$ws.Range("T3").Value = 0.001181759791337041
$ws.Range("U3").Value = 0.001181759791337041
$ws.Range("V3").Value = 0.001181759791337041
$ws.Range("W3").Value = 0.001181759791337041
$ws.Range("X3").Value = 0.001165352380823258
$ws.Range("Y3").Value = 0.001153567294275217

$ws.Range("C4").Value = 0.3339982032775879
$ws.Range("E4").Value = 58.67848997739384
$ws.Range("F4").Value = 0.002158746256690504
$ws.Range("G4").Value = 0.001899229064862406
$ws.Range("H4").Value = 0.001831314705145947
$ws.Range("I4").Value = 0.001654295493918062
$ws.Range("J4").Value = 0.001541817958013691
$ws.Range("K4").Value = 0.001454296567281296
$ws.Range("L4").Value = 0.0013317389100978
$ws.Range("M4").Value = 0.001317833841222117
$ws.Range("N4").Value = 0.001317833841222117
$ws.Range("O4").Value = 0.00128827943957277
$ws.Range("P4").Value = 0.001281367915375082
$ws.Range("Q4").Value = 0.001223858230803504
$ws.Range("R4").Value = 0.001198528213031789
$ws.Range("S4").Value = 0.001198528213031789
$ws.Range("T4").Value = 0.001198528213031789
$ws.Range("U4").Value = 0.001198528213031789
$ws.Range("V4").Value = 0.00118812648764579
$ws.Range("W4").Value = 0.001167655056767268
$ws.Range("X4").Value = 0.001154154023076761
$ws.Range("Y4").Value = 0.001143830213984285

$ws.Range("C5").Value = 0.5180027484893799
$ws.Range("E5").Value = 60.52782793514962
$ws.Range("F5").Value = 0.002194755420862678
$ws.Range("G5").Value = 0.001799809342743957
$ws.Range("H5").Value = 0.001681411725382104
$ws.Range("I5").Value = 0.001681411725382104
$ws.Range("J5").Value = 0.001621076787886049
$ws.Range("K5").Value = 0.001568685401659261
$ws.Range("L5").Value = 0.001458931692714475
$ws.Range("M5").Value = 0.001416287697705302
$ws.Range("N5").Value = 0.001416287697705302
$ws.Range("O5").Value = 0.001315737885109266
$ws.Range("P5").Value = 0.001315737885109266
$ws.Range("Q5").Value = 0.001315737885109266
$ws.Range("R5").Value = 0.001294023936596551
$ws.Range("S5").Value = 0.001257375338803659
$ws.Range("T5").Value = 0.001244273748359075
$ws.Range("U5").Value = 0.001210113468467482
$ws.Range("V5").Value = 0.001197393267106001
$ws.Range("W5").Value = 0.001197393267106001
$ws.Range("X5").Value = 0.001197393267106001
$ws.Range("Y5").Value = 0.001179879686845021

$ws.Range("C6").Value = 0.3619968891143799
$ws.Range("E6").Value = 59.49524954728076
$ws.Range("F6").Value = 0.00224943021628463
$ws.Range("G6").Value = 0.001772159502186653
$ws.Range("H6").Value = 0.001703676119807879
$ws.Range("I6").Value = 0.001490680993088621
$ws.Range("J6").Value = 0.001490680993088621
$ws.Range("K6").Value = 0.001468730539902352
$ws.Range("L6").Value = 0.001468730539902352
$ws.Range("M6").Value = 0.00129916399375485
$ws.Range("N6").Value = 0.00129916399375485
$ws.Range("O6").Value = 0.00129916399375485
$ws.Range("P6").Value = 0.001246147719683289
$ws.Range("Q6").Value = 0.001246147719683289
$ws.Range("R6").Value = 0.001246147719683289
$ws.Range("S6").Value = 0.001220082726477204
$ws.Range("T6").Value = 0.001220082726477204
$ws.Range("U6").Value = 0.001192025421709541
$ws.Range("V6").Value = 0.001159751453163368
$ws.Range("W6").Value = 0.001159751453163368
$ws.Range("X6").Value = 0.001159751453163368
$ws.Range("Y6").Value = 0.001159751453163368

$ws.Range("C7").Value = 0.4499990940093994
$ws.Range("E7").Value = 62.17735807362806
$ws.Range("F7").Value = 0.002146930587310975
$ws.Range("G7").Value = 0.001784917582399589
$ws.Range("H7").Value = 0.00170683299152069
$ws.Range("I7").Value = 0.001655547642938961
$ws.Range("J7").Value = 0.001590170439515019
$ws.Range("K7").Value = 0.001415557859142288
$ws.Range("L7").Value = 0.001397968546080504
$ws.Range("M7").Value = 0.001397968546080504
$ws.Range("N7").Value = 0.001360022177371942
$ws.Range("O7").Value = 0.001293790225855157
$ws.Range("P7").Value = 0.001290332914657192
$ws.Range("Q7").Value = 0.001290332914657192
$ws.Range("R7").Value = 0.001290332914657192
$ws.Range("S7").Value = 0.001290332914657192
$ws.Range("T7").Value = 0.00127104907998969
$ws.Range("U7").Value = 0.00127104907998969
$ws.Range("V7").Value = 0.00124200458369195
$ws.Range("W7").Value = 0.001226590423239226
$ws.Range("X7").Value = 0.001216217441704901
$ws.Range("Y7").Value = 0.001212034270441093

$ws.Range("C8").Value = 0.3400411605834961
$ws.Range("E8").Value = 59.58061462668047
$ws.Range("F8").Value = 0.002207636949791082
$ws.Range("G8").Value = 0.001873654654046744
$ws.Range("H8").Value = 0.001618803621191514
$ws.Range("I8").Value = 0.001618803621191514
$ws.Range("J8").Value = 0.001618803621191514
$ws.Range("K8").Value = 0.001494201536171074
$ws.Range("L8").Value = 0.001494201536171074
$ws.Range("M8").Value = 0.001389272080952055
$ws.Range("N8").Value = 0.001353170658595652
$ws.Range("O8").Value = 0.001318930773223437
$ws.Range("P8").Value = 0.001318930773223437
$ws.Range("Q8").Value = 0.001299827011680392
$ws.Range("R8").Value = 0.001299827011680392
$ws.Range("S8").Value = 0.001238666303088961
$ws.Range("T8").Value = 0.001229334238132207
$ws.Range("U8").Value = 0.001229334238132207
$ws.Range("V8").Value = 0.001161415489798839
$ws.Range("W8").Value = 0.001161415489798839
$ws.Range("X8").Value = 0.001161415489798839
$ws.Range("Y8").Value = 0.001161415489798839

$ws.Range("C9").Value = 0.4089698791503906
$ws.Range("E9").Value = 66.89789480710533
$ws.Range("F9").Value = 0.002201247314085556
$ws.Range("G9").Value = 0.001825589745740979
$ws.Range("H9").Value = 0.001665984951841945
$ws.Range("I9").Value = 0.00158903176056869
$ws.Range("J9").Value = 0.001557972398176961
$ws.Range("K9").Value = 0.001474614463165662
$ws.Range("L9").Value = 0.001418532444851352
$ws.Range("M9").Value = 0.001383823380164133
$ws.Range("N9").Value = 0.001370940766804704
$ws.Range("O9").Value = 0.001351696311884433
$ws.Range("P9").Value = 0.001339925749738076
$ws.Range("Q9").Value = 0.001339925749738076
$ws.Range("R9").Value = 0.001339925749738076
$ws.Range("S9").Value = 0.001339925749738076
$ws.Range("T9").Value = 0.001338625567321422
$ws.Range("U9").Value = 0.001329992490023997
$ws.Range("V9").Value = 0.001326728919210575
$ws.Range("W9").Value = 0.001315189690735935
$ws.Range("X9").Value = 0.00130405253035293
$ws.Range("Y9").Value = 0.00130405253035293

$ws.Range("C10").Value = 0.4319992065429688
$ws.Range("E10").Value = 58.81154154442811
$ws.Range("F10").Value = 0.002235855393236772
$ws.Range("G10").Value = 0.001706908190549372
$ws.Range("H10").Value = 0.001684018799033653
$ws.Range("I10").Value = 0.001623621143254306
$ws.Range("J10").Value = 0.00153288755412732
$ws.Range("K10").Value = 0.00146114651144682
$ws.Range("L10").Value = 0.001452784337469455
$ws.Range("M10").Value = 0.001396883052038582
$ws.Range("N10").Value = 0.001273576252460558
$ws.Range("O10").Value = 0.001273576252460558
$ws.Range("P10").Value = 0.001249162521445244
$ws.Range("Q10").Value = 0.001249162521445244
$ws.Range("R10").Value = 0.001213752816388717
$ws.Range("S10").Value = 0.001213752816388717
$ws.Range("T10").Value = 0.00118912598349265
$ws.Range("U10").Value = 0.001177403052784003
$ws.Range("V10").Value = 0.001165002256058079
$ws.Range("W10").Value = 0.001156940211132454
$ws.Range("X10").Value = 0.001156940211132454
$ws.Range("Y10").Value = 0.001146423811782224

$ws.Range("C11").Value = 0.4550004005432129
$ws.Range("E11").Value = 57.72701545828386
$ws.Range("F11").Value = 0.002217211910674893
$ws.Range("G11").Value = 0.001705433067681969
$ws.Range("H11").Value = 0.001668219450729823
$ws.Range("I11").Value = 0.001630883564591065
$ws.Range("J11").Value = 0.001414394253536275
$ws.Range("K11").Value = 0.001414394253536275
$ws.Range("L11").Value = 0.001305687966391115
$ws.Range("M11").Value = 0.001305687966391115
$ws.Range("N11").Value = 0.001279988611094114
$ws.Range("O11").Value = 0.001270829954691327
$ws.Range("P11").Value = 0.001231056463278251
$ws.Range("Q11").Value = 0.001231056463278251
$ws.Range("R11").Value = 0.00120063610422434
$ws.Range("S11").Value = 0.001187638388142326
$ws.Range("T11").Value = 0.001187638388142326
$ws.Range("U11").Value = 0.001187638388142326
$ws.Range("V11").Value = 0.001159124200348742
$ws.Range("W11").Value = 0.001157437467692809
$ws.Range("X11").Value = 0.00113264965252411
$ws.Range("Y11").Value = 0.001125282952403194

Write-Host "done"